$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 522.1875
$ws.Range("I92").Value = 312.6
$ws.Range("J92").Value = 1270.7142
$ws.Range("K92").Value = 312.6
$ws.Range("L92").Value = 1270.7142
$ws.Range("M92").Value = 935.4
$ws.Range("N92").Value = -3766.7142
$ws.Range("H112").Value = 16034.236
$ws.Range("J112").Value = 17589.66
$ws.Range("L112").Value = 52768.98
$ws.Range("N112").Value = -54984.98
$ws.Range("H129").Value = 19231804
$ws.Range("J129").Value = 1009.6591
$ws.Range("L129").Value = 3028.9773
$ws.Range("N129").Value = -13028.9773
$ws.Range("H135").Value = 1081.7435
$ws.Range("I135").Value = 1108.1621
$ws.Range("K135").Value = 9973.4589
$ws.Range("M135").Value = -7438.4589
$ws.Range("H137").Value = 1365.2307
$ws.Range("I137").Value = 1114.5
$ws.Range("J137").Value = 2511.4285
$ws.Range("K137").Value = 3343.5
$ws.Range("L137").Value = 7534.2855
$ws.Range("M137").Value = -793.5
$ws.Range("N137").Value = -12634.2855
$ws.Range("H138").Value = 3624.6843
$ws.Range("I138").Value = 2156.2307
$ws.Range("J138").Value = 4178.0146
$ws.Range("K138").Value = 6468.6921
$ws.Range("L138").Value = 12534.0438
$ws.Range("M138").Value = -1328.6921
$ws.Range("N138").Value = -22814.0438

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12362401
$ws.Range("I32").Value = 14717824
$ws.Range("J32").Value = 41726.31
$ws.Range("K32").Value = 14717824
$ws.Range("L32").Value = 41726.31
$ws.Range("M32").Value = -14717537
$ws.Range("N32").Value = -42300.31
$ws.Range("H61").Value = 1783.7548
$ws.Range("I61").Value = 1726.8286
$ws.Range("J61").Value = 1894.4445
$ws.Range("K61").Value = 1726.8286
$ws.Range("L61").Value = 1894.4445
$ws.Range("M61").Value = -1514.8286
$ws.Range("N61").Value = -2318.4445
$ws.Range("H74").Value = 1524.3833
$ws.Range("I74").Value = 1631.8372
$ws.Range("J74").Value = 1252.5883
$ws.Range("K74").Value = 1631.8372
$ws.Range("L74").Value = 1252.5883
$ws.Range("M74").Value = -757.8371999999999
$ws.Range("N74").Value = -3000.5883
$ws.Range("H77").Value = 1524.3833
$ws.Range("I77").Value = 1631.8372
$ws.Range("J77").Value = 1252.5883
$ws.Range("K77").Value = 8159.186
$ws.Range("L77").Value = 6262.941499999999
$ws.Range("M77").Value = -3791.186
$ws.Range("N77").Value = -14998.9415
$ws.Range("H97").Value = 1528.1765
$ws.Range("I97").Value = 1462.375
$ws.Range("K97").Value = 1462.375
$ws.Range("M97").Value = -966.375
$ws.Range("H122").Value = 6461.5
$ws.Range("I122").Value = 8315.333000000001
$ws.Range("K122").Value = 24945.999
$ws.Range("M122").Value = -22495.999
$ws.Range("H132").Value = 45460812
$ws.Range("I132").Value = 71429944
$ws.Range("K132").Value = 214289832
$ws.Range("M132").Value = -214287302
$ws.Range("H136").Value = 1783.7548
$ws.Range("I136").Value = 1726.8286
$ws.Range("J136").Value = 1894.4445
$ws.Range("K136").Value = 5180.4858
$ws.Range("L136").Value = 5683.333500000001
$ws.Range("M136").Value = -2630.4858
$ws.Range("N136").Value = -10783.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2911464.5
$ws.Range("I86").Value = 5281
$ws.Range("J86").Value = 7755104
$ws.Range("K86").Value = 5281
$ws.Range("L86").Value = 7755104
$ws.Range("M86").Value = -4158
$ws.Range("N86").Value = -7757350
$ws.Range("H89").Value = 2911464.5
$ws.Range("I89").Value = 5281
$ws.Range("J89").Value = 7755104
$ws.Range("K89").Value = 26405
$ws.Range("L89").Value = 38775520
$ws.Range("M89").Value = -20789
$ws.Range("N89").Value = -38786752
$ws.Range("H94").Value = 751.75
$ws.Range("I94").Value = 658.9091
$ws.Range("J94").Value = 956
$ws.Range("K94").Value = 658.9091
$ws.Range("L94").Value = 956
$ws.Range("M94").Value = -207.9091
$ws.Range("N94").Value = -1858
$ws.Range("H134").Value = 2142701
$ws.Range("I134").Value = 4969.8438
$ws.Range("J134").Value = 5563070.5
$ws.Range("K134").Value = 14909.5314
$ws.Range("L134").Value = 16689211.5
$ws.Range("M134").Value = -12374.5314
$ws.Range("N134").Value = -16694281.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 43479384
$ws.Range("I58").Value = 58824376
$ws.Range("J58").Value = 1916.6666
$ws.Range("K58").Value = 58824376
$ws.Range("L58").Value = 1916.6666
$ws.Range("M58").Value = -58824173
$ws.Range("N58").Value = -2322.6666
$ws.Range("H132").Value = 10421215
$ws.Range("I132").Value = 921.5714
$ws.Range("J132").Value = 30314504
$ws.Range("K132").Value = 2764.7142
$ws.Range("L132").Value = 90943512
$ws.Range("M132").Value = -234.7142000000003
$ws.Range("N132").Value = -90948572
$ws.Range("H134").Value = 1504.75
$ws.Range("I134").Value = 1562.5714
$ws.Range("J134").Value = 1100
$ws.Range("K134").Value = 4687.7142
$ws.Range("L134").Value = 3300
$ws.Range("M134").Value = -2152.7142
$ws.Range("N134").Value = -8370
$ws.Range("H136").Value = 43479384
$ws.Range("I136").Value = 58824376
$ws.Range("J136").Value = 1916.6666
$ws.Range("K136").Value = 176473128
$ws.Range("L136").Value = 5749.9998
$ws.Range("M136").Value = -176470578
$ws.Range("N136").Value = -10849.9998
$ws.Range("H138").Value = 52664.445
$ws.Range("J138").Value = 52664.445
$ws.Range("L138").Value = 52664.445
$ws.Range("N138").Value = -62944.445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 127.89474
$ws.Range("I8").Value = 127.89474
$ws.Range("K8").Value = 383.68422
$ws.Range("M8").Value = -244.68422
$ws.Range("H131").Value = 748.35
$ws.Range("J131").Value = 776.4286
$ws.Range("L131").Value = 2329.2858
$ws.Range("N131").Value = -12409.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6670413.5
$ws.Range("I80").Value = 4419.8
$ws.Range("J80").Value = 20002402
$ws.Range("K80").Value = 4419.8
$ws.Range("L80").Value = 20002402
$ws.Range("M80").Value = -3421.8
$ws.Range("N80").Value = -20004398
$ws.Range("H83").Value = 6670413.5
$ws.Range("I83").Value = 4419.8
$ws.Range("J83").Value = 20002402
$ws.Range("K83").Value = 22099
$ws.Range("L83").Value = 100012010
$ws.Range("M83").Value = -17107
$ws.Range("N83").Value = -100021994
$ws.Range("H113").Value = 3166.6667
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = -2830
$ws.Range("N113").Value = -6590
$ws.Range("H132").Value = 7076.4614
$ws.Range("I132").Value = 1667
$ws.Range("J132").Value = 15731.6
$ws.Range("K132").Value = 5001
$ws.Range("L132").Value = 47194.8
$ws.Range("M132").Value = -2471
$ws.Range("N132").Value = -52254.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 71430136
$ws.Range("I16").Value = 5953399
$ws.Range("J16").Value = 333337060
$ws.Range("K16").Value = 5953399
$ws.Range("L16").Value = 333337060
$ws.Range("M16").Value = -5953229
$ws.Range("N16").Value = -333337400
$ws.Range("H22").Value = 3167506.8
$ws.Range("J22").Value = 5500
$ws.Range("L22").Value = 5500
$ws.Range("N22").Value = -6090
$ws.Range("H27").Value = 3167506.8
$ws.Range("J27").Value = 5500
$ws.Range("L27").Value = 5500
$ws.Range("N27").Value = -5714
$ws.Range("H46").Value = 4167217.5
$ws.Range("I46").Value = 6944910.5
$ws.Range("J46").Value = 678
$ws.Range("K46").Value = 6944910.5
$ws.Range("L46").Value = 678
$ws.Range("M46").Value = -6944722.5
$ws.Range("N46").Value = -1054
$ws.Range("H55").Value = 200020110
$ws.Range("I55").Value = 50051
$ws.Range("J55").Value = 333333500
$ws.Range("K55").Value = 50051
$ws.Range("L55").Value = 333333500
$ws.Range("M55").Value = -49878
$ws.Range("N55").Value = -333333846
$ws.Range("H93").Value = 1365.25
$ws.Range("J93").Value = 1466.3334
$ws.Range("L93").Value = 1466.3334
$ws.Range("N93").Value = -3962.3334
$ws.Range("H100").Value = 3994.625
$ws.Range("I100").Value = 4167.6665
$ws.Range("K100").Value = 4167.6665
$ws.Range("M100").Value = -3626.6665
$ws.Range("H132").Value = 14089336
$ws.Range("I132").Value = 31252172
$ws.Range("J132").Value = 7007.8975
$ws.Range("K132").Value = 93756516
$ws.Range("L132").Value = 21023.6925
$ws.Range("M132").Value = -93753986
$ws.Range("N132").Value = -26083.6925

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("H132").Value = 16969922
$ws.Range("I132").Value = 32292700
$ws.Range("J132").Value = 5417.0713
$ws.Range("K132").Value = 96878100
$ws.Range("L132").Value = 16251.2139
$ws.Range("M132").Value = -96875570
$ws.Range("N132").Value = -21311.2139
$ws.Range("H136").Value = 4899.5293
$ws.Range("I136").Value = 10180.85
$ws.Range("J136").Value = 1492.2258
$ws.Range("K136").Value = 30542.55
$ws.Range("L136").Value = 4476.6774
$ws.Range("M136").Value = -27992.55
$ws.Range("N136").Value = -9576.6774

# Special case: clear M53 on WVR sheet (cell removed from the data in the source diff)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M53").ClearContents()
